$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-33). The value 45188 (2023-09-19) needs to become
# 45189 (2023-09-20) for all of these rows.
for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
